$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 527.5
$ws.Range("I8").Value = 527.5
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1582.5
$ws.Range("L8").Value = 0
$ws.Range("M8").Value = -1443.5
$ws.Range("N8").ClearContents()
$ws.Range("H11").Value = 28.692308
$ws.Range("I11").Value = 28.692308
$ws.Range("K11").Value = 28.692308
$ws.Range("M11").Value = 111.307692
$ws.Range("H31").Value = 2500
$ws.Range("I31").Value = 2500
$ws.Range("K31").Value = 7500
$ws.Range("M31").Value = -7270
$ws.Range("H38").Value = 301.8889
$ws.Range("I38").Value = 214.625
$ws.Range("J38").Value = 1000
$ws.Range("K38").Value = 643.875
$ws.Range("L38").Value = 3000
$ws.Range("M38").Value = -271.875
$ws.Range("N38").Value = -3744
$ws.Range("H40").Value = 6750
$ws.Range("I40").Value = 0
$ws.Range("K40").Value = 0
$ws.Range("M40").ClearContents()
$ws.Range("H51").Value = 9863.272000000001
$ws.Range("I51").Value = 9299.4
$ws.Range("J51").Value = 10333.167
$ws.Range("K51").Value = 9299.4
$ws.Range("L51").Value = 10333.167
$ws.Range("M51").Value = -8815.4
$ws.Range("N51").Value = -11301.167
$ws.Range("H64").Value = 4776.2
$ws.Range("I64").Value = 4814.3335
$ws.Range("K64").Value = 4814.3335
$ws.Range("M64").Value = -4566.3335
$ws.Range("H67").Value = 4776.2
$ws.Range("I67").Value = 4814.3335
$ws.Range("K67").Value = 4814.3335
$ws.Range("M67").Value = -3956.3335
$ws.Range("H113").Value = 2800
$ws.Range("I113").Value = 2800
$ws.Range("K113").Value = 2800
$ws.Range("M113").Value = 454
$ws.Range("H116").Value = 4900
$ws.Range("I116").Value = 4900
$ws.Range("K116").Value = 4900
$ws.Range("M116").Value = -1458
$ws.Range("H132").Value = 792.44446
$ws.Range("I132").Value = 848.125
$ws.Range("J132").Value = 347
$ws.Range("K132").Value = 2544.375
$ws.Range("L132").Value = 1041
$ws.Range("M132").Value = -14.375
$ws.Range("N132").Value = -6101
$ws.Range("H137").Value = 5197.095
$ws.Range("I137").Value = 4952.125
$ws.Range("K137").Value = 14856.375
$ws.Range("M137").Value = -12306.375
$ws.Range("H141").Value = 7191.125
$ws.Range("I141").Value = 7191.125
$ws.Range("K141").Value = 21573.375
$ws.Range("M141").Value = -16393.375

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 585
$ws.Range("J5").Value = 210
$ws.Range("L5").Value = 210
$ws.Range("N5").Value = -434
$ws.Range("H32").Value = 3033.1191
$ws.Range("I32").Value = 3033.1191
$ws.Range("K32").Value = 3033.1191
$ws.Range("M32").Value = -2746.1191
$ws.Range("H45").Value = 5177.7144
$ws.Range("I45").Value = 4707.3335
$ws.Range("K45").Value = 4707.3335
$ws.Range("M45").Value = -4330.3335
$ws.Range("H74").Value = 1801.2
$ws.Range("I74").Value = 1626.5
$ws.Range("K74").Value = 1626.5
$ws.Range("M74").Value = -752.5
$ws.Range("H77").Value = 1801.2
$ws.Range("I77").Value = 1626.5
$ws.Range("K77").Value = 8132.5
$ws.Range("M77").Value = -3764.5
$ws.Range("H132").Value = 5629.857
$ws.Range("I132").Value = 5629.857
$ws.Range("J132").Value = 0
$ws.Range("K132").Value = 16889.571
$ws.Range("L132").Value = 0
$ws.Range("M132").Value = -14359.571
$ws.Range("N132").ClearContents()

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 585
$ws.Range("J4").Value = 210
$ws.Range("L4").Value = 210
$ws.Range("N4").Value = -440
$ws.Range("H86").Value = 3038.1333
$ws.Range("I86").Value = 2969.5
$ws.Range("J86").Value = 3999
$ws.Range("K86").Value = 2969.5
$ws.Range("L86").Value = 3999
$ws.Range("M86").Value = -1846.5
$ws.Range("N86").Value = -6245
$ws.Range("H89").Value = 3038.1333
$ws.Range("I89").Value = 2969.5
$ws.Range("J89").Value = 3999
$ws.Range("K89").Value = 14847.5
$ws.Range("L89").Value = 19995
$ws.Range("M89").Value = -9231.5
$ws.Range("N89").Value = -31227
$ws.Range("H102").Value = 8500
$ws.Range("I102").Value = 8500
$ws.Range("K102").Value = 8500
$ws.Range("M102").Value = -5255
$ws.Range("H134").Value = 6675.6
$ws.Range("I134").Value = 6461.3335
$ws.Range("K134").Value = 19384.0005
$ws.Range("M134").Value = -16849.0005
$ws.Range("H135").Value = 248999.2
$ws.Range("J135").Value = 248999.2
$ws.Range("L135").Value = 248999.2
$ws.Range("N135").Value = -259139.2

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 2179.8
$ws.Range("I31").Value = 1000
$ws.Range("J31").Value = 2474.75
$ws.Range("K31").Value = 1000
$ws.Range("L31").Value = 2474.75
$ws.Range("M31").Value = -705
$ws.Range("N31").Value = -3064.75
$ws.Range("H34").Value = 2179.8
$ws.Range("I34").Value = 1000
$ws.Range("J34").Value = 2474.75
$ws.Range("K34").Value = 1000
$ws.Range("L34").Value = 2474.75
$ws.Range("M34").Value = -798
$ws.Range("N34").Value = -2878.75
$ws.Range("H54").Value = 0
$ws.Range("I54").Value = 0
$ws.Range("J54").Value = 0
$ws.Range("K54").Value = 0
$ws.Range("L54").Value = 0
$ws.Range("M54").ClearContents()
$ws.Range("N54").ClearContents()
$ws.Range("H134").Value = 1637.4242
$ws.Range("I134").Value = 1626.0938
$ws.Range("K134").Value = 4878.2814
$ws.Range("M134").Value = -2343.2814

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 657.5833
$ws.Range("J5").Value = 521.3333
$ws.Range("L5").Value = 1563.9999
$ws.Range("N5").Value = -1787.9999
$ws.Range("H33").Value = 214
$ws.Range("I33").Value = 189
$ws.Range("K33").Value = 1134
$ws.Range("M33").Value = -851
$ws.Range("H55").Value = 8602.333000000001
$ws.Range("I55").Value = 4860
$ws.Range("J55").Value = 10473.5
$ws.Range("K55").Value = 14580
$ws.Range("L55").Value = 31420.5
$ws.Range("M55").Value = -14403
$ws.Range("N55").Value = -31774.5
$ws.Range("H57").Value = 16300
$ws.Range("I57").Value = 1500
$ws.Range("J57").Value = 20000
$ws.Range("K57").Value = 4500
$ws.Range("L57").Value = 60000
$ws.Range("M57").Value = -3941
$ws.Range("N57").Value = -61118
$ws.Range("H122").Value = 666.625
$ws.Range("I122").Value = 489.75
$ws.Range("J122").Value = 843.5
$ws.Range("K122").Value = 4407.75
$ws.Range("L122").Value = 7591.5
$ws.Range("M122").Value = -1957.75
$ws.Range("N122").Value = -12491.5
$ws.Range("H135").Value = 657.5833
$ws.Range("J135").Value = 521.3333
$ws.Range("L135").Value = 4691.9997
$ws.Range("N135").Value = -9761.9997

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 23099.223
$ws.Range("I46").Value = 3000
$ws.Range("J46").Value = 25611.625
$ws.Range("K46").Value = 3000
$ws.Range("L46").Value = 25611.625
$ws.Range("M46").Value = -2844
$ws.Range("N46").Value = -25923.625
$ws.Range("H80").Value = 3333
$ws.Range("I80").Value = 0
$ws.Range("J80").Value = 3333
$ws.Range("K80").Value = 0
$ws.Range("L80").Value = 3333
$ws.Range("M80").ClearContents()
$ws.Range("N80").Value = -5329
$ws.Range("H83").Value = 3333
$ws.Range("I83").Value = 0
$ws.Range("J83").Value = 3333
$ws.Range("K83").Value = 0
$ws.Range("L83").Value = 16665
$ws.Range("M83").ClearContents()
$ws.Range("N83").Value = -26649
$ws.Range("H122").Value = 5149
$ws.Range("I122").Value = 5282.1665
$ws.Range("K122").Value = 15846.4995
$ws.Range("M122").Value = -13396.4995

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1109.7142
$ws.Range("I46").Value = 1347.5
$ws.Range("K46").Value = 1347.5
$ws.Range("M46").Value = -1159.5
$ws.Range("H93").Value = 1980.2727
$ws.Range("I93").Value = 1980.2727
$ws.Range("J93").Value = 0
$ws.Range("K93").Value = 1980.2727
$ws.Range("L93").Value = 0
$ws.Range("M93").Value = -732.2727
$ws.Range("N93").ClearContents()
$ws.Range("H132").Value = 1680.6875
$ws.Range("I132").Value = 1568.826
$ws.Range("K132").Value = 4706.478
$ws.Range("M132").Value = -2176.478

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H2").Value = 50022496
$ws.Range("I2").Value = 50022496
$ws.Range("K2").Value = 50022496
$ws.Range("M2").Value = -50022384
$ws.Range("H81").Value = 3283.611
$ws.Range("I81").Value = 3241.5293
$ws.Range("K81").Value = 6483.0586
$ws.Range("M81").Value = -5422.0586
$ws.Range("H84").Value = 3283.611
$ws.Range("I84").Value = 3241.5293
$ws.Range("K84").Value = 32415.293
$ws.Range("M84").Value = -27111.293
$ws.Range("H96").Value = 2006.3125
$ws.Range("I96").Value = 2254
$ws.Range("J96").Value = 1593.5
$ws.Range("K96").Value = 2254
$ws.Range("L96").Value = 1593.5
$ws.Range("M96").Value = -881
$ws.Range("N96").Value = -4339.5
$ws.Range("H122").Value = 3909.8
$ws.Range("I122").Value = 3909.8
$ws.Range("K122").Value = 11729.4
$ws.Range("M122").Value = -9279.400000000001
$ws.Range("H136").Value = 8646.695
$ws.Range("I136").Value = 8041.619
$ws.Range("K136").Value = 24124.857
$ws.Range("M136").Value = -21574.857
